# Update billing tracking spreadsheet to include billing for the MI royalties.
$wb = $excel.ActiveWorkbook

# --- Non-Collaborative sheet: update the underlying billing amounts ---
$wsNonCollab = $wb.Worksheets.Item("Non-Collaborative")

# Row 2 total formula now includes the new MI royalties amounts (2270 + 480)
$wsNonCollab.Range("B2").Formula = "=2500+2270+480"

# Every other month (rows 3-32) increases from 2150 to 2750
for ($r = 3; $r -le 32; $r++) {
    $wsNonCollab.Cells.Item($r, 2).Value = 2750
}

# The "Total" sheet's B/C columns are formulas that pull from the
# Non-Collaborative and Collaborative sheets, so they recalculate
# automatically once the source data above changes.

# --- Restore the cursor / selected cell that was left on each sheet ---
$wsTotal = $wb.Worksheets.Item("Total")
$wsTotal.Activate()
$wsTotal.Range("E6").Select()

$wsNonCollab.Activate()
$wsNonCollab.Range("B2").Select()

$wsCollab = $wb.Worksheets.Item("Collaborative")
$wsCollab.Activate()
$wsCollab.Range("H42").Select()

# Leave the "Total" sheet as the active tab, matching the original workbook
$wsTotal.Activate()
